$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates pulled from the source diff: (cell, new value, isColumnD)
$updates = @(
    @{ Cell = 'D2'; Value = '57.415.46'; Kind = 'ColD' }
    @{ Cell = 'E2'; Value = '  -0.83%  '; Kind = 'Other' }
    @{ Cell = 'D3'; Value = '3.091.05'; Kind = 'ColD' }
    @{ Cell = 'E3'; Value = '  +0.82%  '; Kind = 'Other' }
    @{ Cell = 'E4'; Value = '  +0.04%  '; Kind = 'Other' }
    @{ Cell = 'D5'; Value = '522.04'; Kind = 'ColD' }
    @{ Cell = 'E5'; Value = '  +0.98%  '; Kind = 'Other' }
    @{ Cell = 'D6'; Value = '140.74'; Kind = 'ColD' }
    @{ Cell = 'E6'; Value = '  -1.21%  '; Kind = 'Other' }
    @{ Cell = 'E7'; Value = '  -0.04%  '; Kind = 'Other' }
    @{ Cell = 'D8'; Value = '3.090.52'; Kind = 'ColD' }
    @{ Cell = 'E8'; Value = '  +0.86%  '; Kind = 'Other' }
    @{ Cell = 'D9'; Value = '0.437'; Kind = 'ColD' }
    @{ Cell = 'E9'; Value = '  +0.16%  '; Kind = 'Other' }
    @{ Cell = 'D10'; Value = '7.18'; Kind = 'ColD' }
    @{ Cell = 'E10'; Value = '  -1.57%  '; Kind = 'Other' }
    @{ Cell = 'D11'; Value = '0.107'; Kind = 'ColD' }
    @{ Cell = 'E11'; Value = '  -0.34%  '; Kind = 'Other' }
    @{ Cell = 'E12'; Value = '  +1.88%  '; Kind = 'Other' }
    @{ Cell = 'D13'; Value = '3.625.04'; Kind = 'ColD' }
    @{ Cell = 'E13'; Value = '  +0.89%  '; Kind = 'Other' }
    @{ Cell = 'E14'; Value = '  +1.09%  '; Kind = 'Other' }
    @{ Cell = 'D15'; Value = '25.93'; Kind = 'ColD' }
    @{ Cell = 'E15'; Value = '  -1.07%  '; Kind = 'Other' }
    @{ Cell = 'E16'; Value = '  -0.69%  '; Kind = 'Other' }
    @{ Cell = 'D17'; Value = '57.488.88'; Kind = 'ColD' }
    @{ Cell = 'E17'; Value = '  -0.69%  '; Kind = 'Other' }
    @{ Cell = 'D18'; Value = '3.098.01'; Kind = 'ColD' }
    @{ Cell = 'E18'; Value = '  +0.98%  '; Kind = 'Other' }
    @{ Cell = 'D19'; Value = '6.10'; Kind = 'ColD' }
    @{ Cell = 'E19'; Value = '  +0.16%  '; Kind = 'Other' }
    @{ Cell = 'D20'; Value = '12.76'; Kind = 'ColD' }
    @{ Cell = 'E20'; Value = '  -0.67%  '; Kind = 'Other' }
    @{ Cell = 'D21'; Value = '8.04'; Kind = 'ColD' }
    @{ Cell = 'E21'; Value = '  -0.77%  '; Kind = 'Other' }
    @{ Cell = 'D22'; Value = '340.38'; Kind = 'ColD' }
    @{ Cell = 'E22'; Value = '  +2.23%  '; Kind = 'Other' }
    @{ Cell = 'E23'; Value = '  +0.11%  '; Kind = 'Other' }
    @{ Cell = 'D24'; Value = '0.511'; Kind = 'ColD' }
    @{ Cell = 'E24'; Value = '  +2.24%  '; Kind = 'Other' }
    @{ Cell = 'D25'; Value = '66.61'; Kind = 'ColD' }
    @{ Cell = 'E25'; Value = '  +1.59%  '; Kind = 'Other' }
    @{ Cell = 'E26'; Value = '  -1.10%  '; Kind = 'Other' }
    @{ Cell = 'D27'; Value = '1.00'; Kind = 'ColD' }
    @{ Cell = 'E27'; Value = '  +0.09%  '; Kind = 'Other' }
    @{ Cell = 'D28'; Value = '0.0₃0906'; Kind = 'ColD' }
    @{ Cell = 'E28'; Value = '  +0.34%  '; Kind = 'Other' }
    @{ Cell = 'D29'; Value = '6.47'; Kind = 'ColD' }
    @{ Cell = 'E30'; Value = '  -0.03%  '; Kind = 'Other' }
    @{ Cell = 'D31'; Value = '7.16'; Kind = 'ColD' }
    @{ Cell = 'E31'; Value = '  -0.76%  '; Kind = 'Other' }
    @{ Cell = 'E32'; Value = '  +1.94%  '; Kind = 'Other' }
    @{ Cell = 'D33'; Value = '20.91'; Kind = 'ColD' }
    @{ Cell = 'E33'; Value = '  +1.02%  '; Kind = 'Other' }
    @{ Cell = 'D34'; Value = '1.18'; Kind = 'ColD' }
    @{ Cell = 'E34'; Value = '  -1.23%  '; Kind = 'Other' }
    @{ Cell = 'D35'; Value = '156.50'; Kind = 'ColD' }
    @{ Cell = 'E35'; Value = '  +1.01%  '; Kind = 'Other' }
    @{ Cell = 'D36'; Value = '4.62'; Kind = 'ColD' }
    @{ Cell = 'E36'; Value = '  +1.80%  '; Kind = 'Other' }
    @{ Cell = 'E37'; Value = '  +1.07%  '; Kind = 'Other' }
    @{ Cell = 'D38'; Value = '27.08'; Kind = 'ColD' }
    @{ Cell = 'E38'; Value = '  +0.72%  '; Kind = 'Other' }
    @{ Cell = 'D39'; Value = '1.27'; Kind = 'ColD' }
    @{ Cell = 'E39'; Value = '  +0.42%  '; Kind = 'Other' }
    @{ Cell = 'D40'; Value = '0.0656'; Kind = 'ColD' }
    @{ Cell = 'E40'; Value = '  -2.98%  '; Kind = 'Other' }
    @{ Cell = 'B41'; Value = 'Stacks'; Kind = 'Other' }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'; Kind = 'Other' }
    @{ Cell = 'D41'; Value = '1.52'; Kind = 'ColD' }
    @{ Cell = 'E41'; Value = '  +11.15%  '; Kind = 'Other' }
    @{ Cell = 'D42'; Value = '3.135.83'; Kind = 'ColD' }
    @{ Cell = 'E42'; Value = '  +0.91%  '; Kind = 'Other' }
    @{ Cell = 'B43'; Value = 'Filecoin'; Kind = 'Other' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; Kind = 'Other' }
    @{ Cell = 'D43'; Value = '3.93'; Kind = 'ColD' }
    @{ Cell = 'E43'; Value = '  +0.32%  '; Kind = 'Other' }
    @{ Cell = 'B44'; Value = 'Mantle'; Kind = 'Other' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; Kind = 'Other' }
    @{ Cell = 'D44'; Value = '0.684'; Kind = 'ColD' }
    @{ Cell = 'E44'; Value = '  +4.44%  '; Kind = 'Other' }
    @{ Cell = 'D45'; Value = '36.65'; Kind = 'ColD' }
    @{ Cell = 'E45'; Value = '  +0.51%  '; Kind = 'Other' }
    @{ Cell = 'E46'; Value = '  +0.05%  '; Kind = 'Other' }
    @{ Cell = 'D47'; Value = '2.312.10'; Kind = 'ColD' }
    @{ Cell = 'E47'; Value = '  +2.09%  '; Kind = 'Other' }
    @{ Cell = 'E48'; Value = '  -0.15%  '; Kind = 'Other' }
    @{ Cell = 'D49'; Value = '0.972'; Kind = 'ColD' }
    @{ Cell = 'E49'; Value = '  +2.68%  '; Kind = 'Other' }
    @{ Cell = 'D50'; Value = '20.64'; Kind = 'ColD' }
    @{ Cell = 'E50'; Value = '  -0.46%  '; Kind = 'Other' }
    @{ Cell = 'D51'; Value = '6.00'; Kind = 'ColD' }
    @{ Cell = 'E51'; Value = '  +1.28%  '; Kind = 'Other' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Kind -eq "ColD") {
        # Column D holds price strings that look numeric (e.g. "1.00", "6.10").
        # Force text so Excel keeps the exact string instead of coercing to a
        # number and dropping trailing zeros / the thousands-dot formatting,
        # then drop back to the Normal style so no stray number format sticks.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}

Write-Output "Applied $($updates.Count) cell updates"
